$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author was debugging the "koperasi" import path, so the two extra
# sample products (rows 3 and 4) are thrown away entirely ...
$ws.Rows("3:4").Delete()

# ... and the remaining product row is overwritten with throwaway debug
# values: a plain (non-string) number in A2, -1 in B2, and 1.25 in C2 with
# its number format explicitly (re)applied.
$ws.Range("A2").Value = 3.14159
$ws.Range("B2").Value = -1
$ws.Range("C2").NumberFormat = "General"
$ws.Range("C2").Value = 1.25

# Keep the saved selection in sync with where Excel last left the cursor.
$ws.Range("C2").Select() | Out-Null
